# "updates to PS layout"
# Parts are back in stock / no longer "missing" -- mark their status as "ok".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value  = "ok"
$ws.Range("F19").Value = "ok"
$ws.Range("F24").Value = "ok"
$ws.Range("F53").Value = "ok"
$ws.Range("F57").Value = "ok"

# Re-position the view: scroll the grid down a bit and move the selection
# to the new bottom of the list.
$ws.Activate()
$ws.Range("A47").Select()
$ws.Range("F72").Select()

try {
    $excel.ActiveWindow.ScrollRow = 47
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

try {
    $excel.ActiveWindow.Left   = 0
    $excel.ActiveWindow.Top    = 0
    $excel.ActiveWindow.Width  = 25600
    $excel.ActiveWindow.Height = 15460
} catch {}

"done"
